$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared/rich-text header updates ---
$ws.Range("A8").Value = "Volume 31   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/15/2024  Through  1/21/2024"

# --- Weekly crime statistics table (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = -50
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = -16.666666666666
$ws.Range("I14").Value = 4
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = -20
$ws.Range("L14").Value = 33.333333333333
$ws.Range("M14").Value = -20
$ws.Range("N14").Value = -87.096774193548

# Row 15
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 10
$ws.Range("G15").Value = 27
$ws.Range("H15").Value = -62.962962962963
$ws.Range("I15").Value = 7
$ws.Range("J15").Value = 21
$ws.Range("K15").Value = -66.666666666666
$ws.Range("L15").Value = -41.666666666666
$ws.Range("M15").Value = -36.363636363636
$ws.Range("N15").Value = -77.419354838709

# Row 16
$ws.Range("C16").Value = 40
$ws.Range("D16").Value = 39
$ws.Range("E16").Value = 2.564102564102
$ws.Range("F16").Value = 179
$ws.Range("G16").Value = 158
$ws.Range("H16").Value = 13.291139240506
$ws.Range("I16").Value = 129
$ws.Range("J16").Value = 131
$ws.Range("K16").Value = -1.526717557251
$ws.Range("L16").Value = 3.2
$ws.Range("M16").Value = -38.277511961722
$ws.Range("N16").Value = -87.189672293942

# Row 17
$ws.Range("C17").Value = 58
$ws.Range("D17").Value = 69
$ws.Range("E17").Value = -15.942028985507
$ws.Range("F17").Value = 269
$ws.Range("G17").Value = 273
$ws.Range("H17").Value = -1.465201465201
$ws.Range("I17").Value = 202
$ws.Range("J17").Value = 215
$ws.Range("K17").Value = -6.046511627906
$ws.Range("L17").Value = -0.492610837438
$ws.Range("M17").Value = 23.170731707317
$ws.Range("N17").Value = -55.604395604395

# Row 18
$ws.Range("C18").Value = 39
$ws.Range("D18").Value = 55
$ws.Range("E18").Value = -29.090909090909
$ws.Range("F18").Value = 154
$ws.Range("G18").Value = 215
$ws.Range("H18").Value = -28.372093023255
$ws.Range("I18").Value = 114
$ws.Range("J18").Value = 165
$ws.Range("K18").Value = -30.90909090909
$ws.Range("L18").Value = -27.848101265822
$ws.Range("M18").Value = -28.75
$ws.Range("N18").Value = -82.407407407407

# Row 19
$ws.Range("C19").Value = 75
$ws.Range("D19").Value = 123
$ws.Range("E19").Value = -39.024390243902
$ws.Range("F19").Value = 329
$ws.Range("G19").Value = 428
$ws.Range("H19").Value = -23.130841121495
$ws.Range("I19").Value = 242
$ws.Range("J19").Value = 325
$ws.Range("K19").Value = -25.538461538461
$ws.Range("L19").Value = -19.601328903654
$ws.Range("M19").Value = 30.81081081081
$ws.Range("N19").Value = -30.259365994236

# Row 20
$ws.Range("C20").Value = 35
$ws.Range("D20").Value = 32
$ws.Range("E20").Value = 9.375
$ws.Range("F20").Value = 135
$ws.Range("G20").Value = 117
$ws.Range("H20").Value = 15.384615384615
$ws.Range("I20").Value = 109
$ws.Range("J20").Value = 90
$ws.Range("K20").Value = 21.111111111111
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 39.743589743589
$ws.Range("N20").Value = -81.803005008347

# Row 21
$ws.Range("C21").Value = 251
$ws.Range("D21").Value = 326
$ws.Range("E21").Value = -23.006134969325
$ws.Range("F21").Value = 1081
$ws.Range("G21").Value = 1224
$ws.Range("H21").Value = -11.683006535947
$ws.Range("I21").Value = 807
$ws.Range("J21").Value = 952
$ws.Range("K21").Value = -15.231092436974
$ws.Range("L21").Value = -11.416026344676
$ws.Range("M21").Value = -0.615763546798
$ws.Range("N21").Value = -74.118024374599

# Row 22
$ws.Range("C22").Value = 6
$ws.Range("D22").Value = 8
$ws.Range("E22").Value = -25
$ws.Range("F22").Value = 26
$ws.Range("G22").Value = 28
$ws.Range("H22").Value = -7.142857142857
$ws.Range("I22").Value = 20
$ws.Range("J22").Value = 21
$ws.Range("K22").Value = -4.761904761904
$ws.Range("L22").Value = -9.090909090909
$ws.Range("M22").Value = 0

# Row 23
$ws.Range("C23").Value = 29
$ws.Range("D23").Value = 32
$ws.Range("E23").Value = -9.375
$ws.Range("F23").Value = 104
$ws.Range("G23").Value = 118
$ws.Range("H23").Value = -11.864406779661
$ws.Range("I23").Value = 78
$ws.Range("J23").Value = 92
$ws.Range("K23").Value = -15.217391304347
$ws.Range("L23").Value = -13.333333333333
$ws.Range("M23").Value = 62.5

# Row 24
$ws.Range("C24").Value = 173
$ws.Range("D24").Value = 228
$ws.Range("E24").Value = -24.122807017543
$ws.Range("F24").Value = 811
$ws.Range("G24").Value = 879
$ws.Range("H24").Value = -7.736063708759
$ws.Range("I24").Value = 575
$ws.Range("J24").Value = 666
$ws.Range("K24").Value = -13.663663663663
$ws.Range("L24").Value = -6.351791530944
$ws.Range("M24").Value = 0.877192982456

# Row 25
$ws.Range("C25").Value = 102
$ws.Range("D25").Value = 122
$ws.Range("E25").Value = -16.39344262295
$ws.Range("F25").Value = 425
$ws.Range("G25").Value = 447
$ws.Range("H25").Value = -4.921700223713
$ws.Range("I25").Value = 313
$ws.Range("J25").Value = 334
$ws.Range("K25").Value = -6.2874251497
$ws.Range("L25").Value = 3.986710963455
$ws.Range("M25").Value = -24.213075060532

# Row 26
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 14.285714285714
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = -27.272727272727
$ws.Range("I26").Value = 20
$ws.Range("J26").Value = 26
$ws.Range("K26").Value = -23.076923076923
$ws.Range("L26").Value = -4.761904761904

# Row 27
$ws.Range("C27").Value = 11
$ws.Range("D27").Value = 13
$ws.Range("E27").Value = -15.384615384615
$ws.Range("F27").Value = 38
$ws.Range("G27").Value = 46
$ws.Range("H27").Value = -17.391304347826
$ws.Range("I27").Value = 26
$ws.Range("J27").Value = 39
$ws.Range("K27").Value = -33.333333333333
$ws.Range("L27").Value = 13.043478260869

# Row 28
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = 25
$ws.Range("F28").Value = 14
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = 40
$ws.Range("I28").Value = 11
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = 37.5
$ws.Range("L28").Value = -8.333333333333
$ws.Range("M28").Value = -42.105263157894
$ws.Range("N28").Value = -90.677966101694

# Row 29
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 12
$ws.Range("G29").Value = 10
$ws.Range("H29").Value = 20
$ws.Range("I29").Value = 10
$ws.Range("J29").Value = 8
$ws.Range("K29").Value = 25
$ws.Range("L29").Value = -9.090909090909
$ws.Range("M29").Value = -33.333333333333
$ws.Range("N29").Value = -90.825688073394

# Row 30
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 5
$ws.Range("E30").Value = -80
$ws.Range("F30").Value = 8
$ws.Range("G30").Value = 7
$ws.Range("H30").Value = 14.285714285714
$ws.Range("I30").Value = 7
$ws.Range("J30").Value = 7
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 40
